$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.159670948982239
$ws.Range("B1").Value = 2.370934724807739
$ws.Range("D1").Value = 2.396783351898193
$ws.Range("E1").Value = 1.221003770828247
